$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TDSP")

# Rows 2-101: the value/value_t columns (B,C) are being cleared out (stale
# duplicate data removed upstream). ClearContents is the faithful Excel-COM
# equivalent of "no longer populated" for these historical rows.
$ws.Range("B2:C101").ClearContents()

# Rows 148-150 and 161-183: refreshed data values (the source series was
# re-pulled and slightly revised numbers flowed in for the recent history).
$ws.Range("B148:C148").Value = 11.769828
$ws.Range("B149:C149").Value = 11.865791
$ws.Range("B150:C150").Value = 11.723736

$ws.Range("B161:C161").Value = 11.727541
$ws.Range("B162:C162").Value = 11.591065
$ws.Range("B163:C163").Value = 9.739786
$ws.Range("B164:C164").Value = 10.058434
$ws.Range("B165:C165").Value = 10.389279
$ws.Range("B166:C166").Value = 9.051398000000001
$ws.Range("B167:C167").Value = 9.840935999999999
$ws.Range("B168:C168").Value = 10.009535
$ws.Range("B169:C169").Value = 10.228733
$ws.Range("B170:C170").Value = 10.47232
$ws.Range("B171:C171").Value = 10.684637
$ws.Range("B172:C172").Value = 10.567735
$ws.Range("B173:C173").Value = 10.736833
$ws.Range("B174:C174").Value = 10.56399
$ws.Range("B175:C175").Value = 10.576866
$ws.Range("B176:C176").Value = 10.747704
$ws.Range("B177:C177").Value = 11.095984
$ws.Range("B178:C178").Value = 11.058721
$ws.Range("B179:C179").Value = 11.019084
$ws.Range("B180:C180").Value = 11.138535
$ws.Range("B181:C181").Value = 11.12225
$ws.Range("B182:C182").Value = 11.10531
$ws.Range("B183:C183").Value = 11.123837

# New row 184: the latest quarterly observation (2025-07-01, serial 45839).
$ws.Range("A184").Value = 45839
$ws.Range("A184").NumberFormat = $ws.Range("A183").NumberFormat
$ws.Range("B184:C184").Value = 11.256338
